# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interested count) / "最低票价" (lowest price) /
# cover-image values to the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  F = 11434 },
    @{ Row = 3;  F = 10817 },
    @{ Row = 5;  F = 4 },
    @{ Row = 6;  F = 983 },
    @{ Row = 8;  F = 53 },
    @{ Row = 10; F = 36 },
    @{ Row = 11; F = 10586 },
    @{ Row = 12; F = 4085 },
    @{ Row = 13; F = 7; G = 45; I = "//i2.hdslb.com/bfs/openplatform/202401/lVVrbSra1706508320671.jpeg" },
    @{ Row = 14; F = 2452 },
    @{ Row = 15; F = 45 },
    @{ Row = 16; F = 28 },
    @{ Row = 17; F = 106 },
    @{ Row = 18; F = 417 },
    @{ Row = 19; F = 11093 },
    @{ Row = 20; F = 10851 }
)

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Range("F" + $u.Row).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Range("G" + $u.Row).Value = $u.G
        }
        if ($u.ContainsKey("I")) {
            $ws.Range("I" + $u.Row).Value = $u.I
        }
    }
}
